$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price cells are stored as literal text even when the string looks numeric
# (e.g. "1.000", "305.31"). Force text format while writing those values so Excel
# does not silently convert them to numbers and strip significant trailing zeros,
# then restore the default "Normal" style so no stray number format lingers on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9994"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3760"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "53.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3611"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.266"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08147"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.616"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.365"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001245"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.543"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.412"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.272"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.402"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.825"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9545"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02768"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07391"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2512"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.126"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08769"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7104"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6537"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.331"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9982"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.011"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "133.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07974"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.199"
$ws.Range("D51").Style = "Normal"

# Remaining updated cells (Volume(1h) percentages, and Price cells that are safely
# non-numeric strings already, e.g. "23.238.48") can be assigned directly.
$ws.Range("D2").Value = "23.238.48"
$ws.Range("D3").Value = "1.608.48"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  +4.49%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("E13").Value = "  +3.18%  "
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "1.604.09"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "23.243.10"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("E25").Value = "  +10.34%  "
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").Value = "1.785.00"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  +2.49%  "
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("E44").Value = "  +5.95%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  +2.74%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  -1.19%  "
